# "Work on budget and time plan"
#
# Adds a "Semester 1" label under the title, revises several of the
# estimated-hours figures in the A: SALARY table (which ripple through the
# dependent Cost/total/quote formulas automatically), updates the date
# columns' number format to the built-in d-mmm-yy format, and leaves the
# final selection on B23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Semester 1" label under the report title (row 2 was previously blank).
$ws.Range("B2").Value = "Semester 1"

# Revised hour estimates for rows 16-20 of the salary table; the Cost column
# (E16:E20) is a shared formula (=D*$D$9) and the grand totals / quote figures
# below (D24, E24, C31, C32, C33) recalculate automatically.
$ws.Range("D16").Value = 13
$ws.Range("D17").Value = 33
$ws.Range("D18").Value = 20
$ws.Range("D19").Value = 40
$ws.Range("D20").Value = 20

# The blank "Start date" cells switch from the custom dd-mmm-yy format to
# Excel's built-in d-mmm-yy date format (built-in numFmtId 15).
$ws.Range("C13").NumberFormat = "d-mmm-yy"
$ws.Range("C14").NumberFormat = "d-mmm-yy"
$ws.Range("C15:C22").NumberFormat = "d-mmm-yy"
$ws.Range("C24").NumberFormat = "d-mmm-yy"

# Restore the view: scrolled near the top of the salary table, with B23 as
# the active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("B23").Select()
